$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: drop the " [%]" suffix from each metric name
$ws.Range("E4").Value = "lemmatizer accuracy"
$ws.Range("F4").Value = "ner precision"
$ws.Range("G4").Value = "ner recall"
$ws.Range("H4").Value = "ner F1 score"

# Update the active selection to H4 (matches saved selection in the file)
$ws.Range("H4").Select()
